$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, shifting existing rows 134-211 down to 135-212
$ws.Rows("134:134").Insert()

# Populate the newly inserted row 134 with the new "Sweet Heart" cherry record
$ws.Range("A134").Value = 10
$ws.Range("B134").Value = "Vega Modelo de Temuco"
$ws.Range("C134").Value = "La Araucanía"
$ws.Range("D134").Value = 44572
$ws.Range("E134").Value = 9
$ws.Range("F134").Value = "Fruta"
$ws.Range("G134").Value = 100103
$ws.Range("H134").Value = "Frutos de hueso (carozo)"
$ws.Range("I134").Value = 100103001
$ws.Range("J134").Value = "Cereza"
$ws.Range("K134").Value = "Sweet Heart"
$ws.Range("L134").Value = "Primera"
$ws.Range("M134").Value = 800
$ws.Range("N134").Value = 6000
$ws.Range("O134").Value = 6000
$ws.Range("P134").Value = 6000
$ws.Range("Q134").Value = "$/bandeja 6 kilos"
$ws.Range("R134").Value = "Región del Maule"
$ws.Range("S134").Value = 1000
$ws.Range("T134").Value = 6
